$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append (rows 20..28), matching columns:
# A: index, B: Model name, C: Score, D: Acc, E: Dataset, F: Params,
# G: Compute time, H: Model, I: n_gram, J: feat_type, K: lemmatized,
# L: lang, M: standardized, N: rm_accents, O: feat_select, P: n_feat
$rows = @(
    @{ A=18; B="ComplementNB"; C=0.7397804972804973; D=0.8303198887343533; E="L - TF - NLe - 2G - F2000"; F="{}"; G=0.1761283874511719; H="ComplementNB()"; I="(1, 2)"; J="TF";    K=$false; L=$true;  M=$false; N=$false; O="F_CL"; P=2000 },
    @{ A=19; B="ComplementNB"; C=0.7934702797202797; D=0.8831710709318498; E="L - C - NLe - 2G - F2000 - No Acc";        F="{}"; G=0.2726867198944092; H="ComplementNB()"; I="(1, 2)"; J="Count"; K=$false; L=$true;  M=$false; N=$true;  O="F_CL"; P=2000 },
    @{ A=20; B="ComplementNB"; C=0.7932915695415695; D=0.8831710709318498; E="L - C - NLe - 2G - F2000 - Punc w space"; F="{}"; G=0.1486952304840088; H="ComplementNB()"; I="(1, 2)"; J="Count"; K=$false; L=$true;  M=$false; N=$false; O="F_CL"; P=2000 },
    @{ A=21; B="ComplementNB"; C=0.6340540015540015; D=1;                  E="NL - B - NLe - 2G";                      F="{}"; G=0.4496505260467529; H="ComplementNB()"; I="(1, 2)"; J="Bin";   K=$false; L=$false; M=$false; N=$false; O="";     P=1    },
    @{ A=22; B="ComplementNB"; C=0.6576991064491065; D=1;                  E="NL - C - NLe - 2G";                      F="{}"; G=0.4078330993652344; H="ComplementNB()"; I="(1, 2)"; J="Count"; K=$false; L=$false; M=$false; N=$false; O="";     P=1    },
    @{ A=23; B="ComplementNB"; C=0.8034741647241648; D=0.885952712100139;  E="Big test";                               F="{}"; G=0.1550121307373047; H="ComplementNB()"; I="(1, 2)"; J="Count"; K=$false; L=$true;  M=$false; N=$true;  O="F_CL"; P=2000 },
    @{ A=24; B="ComplementNB"; C=0.8061120823620824; D=0.8929068150208623; E="Big test- BIN";                          F="{}"; G=0.1614856719970703; H="ComplementNB()"; I="(1, 2)"; J="Bin";   K=$false; L=$true;  M=$false; N=$true;  O="F_CL"; P=2000 },
    @{ A=25; B="ComplementNB"; C=0.8079254079254079; D=0.8984700973574409; E="Big test- BIN - NL";                     F="{}"; G=0.1612675189971924; H="ComplementNB()"; I="(1, 2)"; J="Bin";   K=$false; L=$false; M=$false; N=$true;  O="F_CL"; P=2000 },
    @{ A=26; B="ComplementNB"; C=0.8063908313908315; D=0.8901251738525731; E="Big test- Count - NL";                   F="{}"; G=0.1819298267364502; H="ComplementNB()"; I="(1, 2)"; J="Count"; K=$false; L=$false; M=$false; N=$true;  O="F_CL"; P=2000 }
)

$styleSource = $ws.Range("A19")
$styleSource.Copy()

$startRow = 20
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $aCell = $ws.Cells.Item($r, 1)
    $aCell.Value = $data.A
    $aCell.PasteSpecial(-4122)

    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 3).Value = $data.C
    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 5).Value = $data.E
    $ws.Cells.Item($r, 6).Value = $data.F
    $ws.Cells.Item($r, 7).Value = $data.G
    $ws.Cells.Item($r, 8).Value = $data.H
    $ws.Cells.Item($r, 9).Value = $data.I
    $ws.Cells.Item($r, 10).Value = $data.J
    $ws.Cells.Item($r, 11).Value = $data.K
    $ws.Cells.Item($r, 12).Value = $data.L
    $ws.Cells.Item($r, 13).Value = $data.M
    $ws.Cells.Item($r, 14).Value = $data.N
    $ws.Cells.Item($r, 15).Value = $data.O
    $ws.Cells.Item($r, 16).Value = $data.P
}
